$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JdT-TPI_LRD")

# Insert two new rows by duplicating the formatting of the last existing
# data row (row 46) so the date-format / wrap-text styles carry over
# exactly, then overwrite the values for the two new entries.
$ws.Rows("46:46").Copy()
$ws.Rows("47:47").Insert(-4121)
$ws.Rows("46:46").Copy()
$ws.Rows("48:48").Insert(-4121)
$excel.CutCopyMode = 0

$ws.Range("A47").Value = 44693
$ws.Range("B47").Value = "Réalisation"
$ws.Range("C47").Value = 0.75
$ws.Range("D47").Value = "Commencement de la fonction de login"

$ws.Range("A48").Value = 44694
$ws.Range("B48").Value = "Réalisation"
$ws.Range("C48").Value = 0.75
$ws.Range("D48").Value = "Finition de la fonction de login"

$ws.ListObjects.Item("Tableau1").Resize($ws.Range("A1:F48"))

$ws.Range("C47").Select()
